# week 10 sum 2022 updates
# Append new matchup rows (1091-1109) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(6,1,7,2),
    @(3,3,3,0),
    @(4,3,3,0),
    @(4,3,4,0),
    @(5,2,5,1),
    @(6,1,6,2),
    @(2,2,2,1),
    @(6,2,6,0),
    @(4,2,5,1),
    @(3,3,4,0),
    @(5,0,5,2),
    @(5,2,7,0),
    @(4,3,3,0),
    @(5,2,5,0),
    @(2,2,3,0),
    @(7,1,6,2),
    @(3,2,3,1),
    @(5,3,4,0),
    @(2,2,3,0)
)

$startRow = 1091
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}

$excel.ActiveWindow.ScrollRow = 1088
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G1108").Select()
